# Regenerate the localization-status report for the latest handoff run.
#
# This mirrors the CI tool's report-generation step: a fresh handoff pass
# bumped the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the 1d9d794f... file (and the other rows that happened to
# share those exact timestamps from the same batch run), and the six
# rows that were missing a Priority got stamped with the "ht" handoff type.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn     = $wb.Sheets.Item("zh-cn")
$dede     = $wb.Sheets.Item("de-de")

# Rows (in each localized-language table) whose "Priority" column (E) needs
# to be stamped with the handoff type "ht" instead of being blank.
$priorityRows = @(7, 9, 11, 12, 13, 14)

foreach ($row in $priorityRows) {
    $zhcn.Range("E$row").Value = "ht"
    $dede.Range("E$row").Value = "ht"
}

# Rows on the Overview sheet whose "Latest HO Xliff Generate Date" (G) carried
# the old "2016-08-27 12:19:38" timestamp now carry the refreshed one.
$overviewDateRows = @(7, 9, 11, 12, 13, 14)
foreach ($row in $overviewDateRows) {
    $overview.Range("G$row").Value = "2016-08-27 12:19:53"
}

# The de-de table's "Latest Handoff Datetime" (H) column mirrored that same
# old Overview timestamp on these rows, so it gets refreshed too.
$dedeDateRows = @(7, 9, 11, 12, 13, 14)
foreach ($row in $dedeDateRows) {
    $dede.Range("H$row").Value = "2016-08-27 12:19:53"
}

# The zh-cn table's "Latest Handoff Datetime" (H) column had its own,
# slightly earlier timestamp ("2016-08-27 12:19:34") on the same rows.
$zhcnDateRows = @(7, 9, 11, 12, 13, 14)
foreach ($row in $zhcnDateRows) {
    $zhcn.Range("H$row").Value = "2016-08-27 12:19:49"
}
